$d = $word.ActiveDocument

# 1. Title paragraph: switch from direct bold/size formatting to the "Titolo1" (Heading 1) style.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Style = "Titolo1"

# 2. Split the run " quella serie di note ordinate in senso ascendente o discendente da un suono..."
#    into three runs by inserting "che va " before "da un suono...".
$rng = $d.Content
$rng.Find.Execute("da un suono sino a quello avete lo stesso nome", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $rng.Start
$ins = $d.Range($insertPoint, $insertPoint)
$ins.InsertBefore("che va ")
# Force a run boundary between "discendente " and "che va " by toggling a
# character formatting property on the newly inserted text and back again.
$newRunRange = $d.Range($insertPoint, $insertPoint + 7)
$newRunRange.Font.Bold = 1
$newRunRange.Font.Bold = 0

# 3. Move the "_GoBack" bookmark from the end of the "nota fondamentale" paragraph
#    to right after "...la seguente serie di intervalli: " in the "scala minore" paragraph.
$bmRng = $d.Content
$bmRng.Find.Execute("quella scala la cui sequenza di note ordinate in senso ascendente dà la seguente serie di intervalli: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPoint = $bmRng.End
$bmTarget = $d.Range($bmPoint, $bmPoint)
$d.Bookmarks.Add("_GoBack", $bmTarget)
